# Add team record (Wins / Losses / Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (currently rows 1..43, header + 42 players).
$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
if ($lastRow -lt 2) { $lastRow = 43 }

# Copy the header style from the existing last header cell (AC1) onto the
# three new header cells so they match the bold/bordered/centered look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill every data row (2 through the last row) with the team's record.
$ws.Range("AD2:AD" + $lastRow).Value = 88
$ws.Range("AE2:AE" + $lastRow).Value = 74
$ws.Range("AF2:AF" + $lastRow).Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
